# LoginData.xlsx - Added RPA tests for Excel Read, Write in Browser tests + bonus keywords
#
# - Removes the now-unused "Data2" worksheet
# - Adds a "Result" column workflow on the "login" sheet (new C2 value)
# - Rebuilds "Data1" with the new read/write RPA test-data layout

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1. Drop the "Data2" worksheet entirely
# ------------------------------------------------------------------
[void]$wb.Worksheets("Data2").Delete()

# ------------------------------------------------------------------
# 2. "login" sheet: add a Result column with an extra row for the
#    invalid-credentials case
# ------------------------------------------------------------------
$login = $wb.Worksheets("login")

$login.Range("C1").Value = "Result"
$login.Range("C2").Value = "Invalid credentials"

$login.Columns("C").ColumnWidth = 17

[void]$login.Range("C2").Select()

# ------------------------------------------------------------------
# 3. "Data1" sheet: replace the old Test1/Test2 list with the new
#    RPA read/write test-data grid
# ------------------------------------------------------------------
$data1 = $wb.Worksheets("Data1")

$data1.Cells.ClearContents()

$data1.Range("D1").Value = "A1"
$data1.Range("G1").Value = "B"

$data1.Range("D2").Value = "A2"
$data1.Range("G2").Value = "B"

$data1.Range("D7").Value = "A3"
$data1.Range("G7").Value = "B"

$data1.Range("D8").Value = "A7"
$data1.Range("G8").Value = "B"

[void]$data1.Activate()
[void]$data1.Range("D6:D7").Select()
